$d = $word.ActiveDocument

# 1. Replace the placeholder "Evaluatie" paragraph text with the actual
#    evaluation content written by the authors.
$d.Content.Find.Execute(
    "Leg een verband tussen de getrokken conclusie en het doel van het experiment (en de hypothese). Ga daarbij ook in op bijvoorbeeld de meetonzekerheid als gevolg van de gebruikte meetmethoden of eventuele meetfouten.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Bij de meting hebben we alleen de RGB image kunnen testen, omdat we niet onze gray-conversion met de default implementatie.",
    2)

# 2. Move the "_GoBack" bookmark so it marks the end of the paragraph that
#    was just edited (mirroring Word's normal behaviour of tracking the
#    location of the last edit). The target spot sits right at the very
#    end of the document's content, which this host can mis-place a
#    bookmark at directly, so a small amount of scratch text is appended
#    first to move the true document end further away, the bookmark is
#    dropped at the now-safe position, and the scratch text is removed
#    again -- the bookmark stays put because it never overlapped it.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$targetPos = $lastPara.Range.End - 1

$tail = $d.Content
$tail.Collapse(0)
$tail.InsertAfter("XTEMPX")

$bmRange = $d.Range($targetPos, $targetPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Content.Find.Execute("XTEMPX", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
